$d = $word.ActiveDocument

# The sentence "Makes HTTP requests to the API for data." needs to lose its
# trailing period and gain a new, identically-formatted run describing the
# axios API: "Makes HTTP requests to the API for data" + " (axios API)."
$finder = $d.Content
$found = $finder.Find.Execute("Makes HTTP requests to the API for data.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-materialize the Find hit as a plain Range: Range objects returned
    # directly from Find.Execute append rather than replace when used with
    # InsertXML, so rebuild an equivalent Range from Start/End first.
    $target = $d.Range($finder.Start, $finder.End)

    $runPr = '<w:rPr>' +
      '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
      '<w:kern w:val="0"/>' +
      '<w:sz w:val="24"/>' +
      '<w:szCs w:val="24"/>' +
      '<w:lang w:eastAsia="en-CA"/>' +
      '<w14:ligatures w14:val="none"/>' +
      '</w:rPr>'

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
      'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
      '<w:r w:rsidRPr="009973A7">' + $runPr +
      '<w:t>Makes HTTP requests to the API for data</w:t>' +
      '</w:r>' +
      '<w:r>' + $runPr +
      '<w:t xml:space="preserve"> (axios API).</w:t>' +
      '</w:r>' +
      '</w:p>'

    try {
        $target.InsertXML($xml) | Out-Null
    } catch {
        # Fall back to a plain-text replacement if InsertXML is unavailable;
        # the visible wording still ends up correct even though the run
        # split from the original diff would be lost in that case.
        $target.Text = "Makes HTTP requests to the API for data (axios API)."
    }
}
